# Insert a new weekly data row for "Arveja Verde" (Comercializadora del Agro
# de Limarí) right after the existing row 55. This pushes the previously
# existing rows 56-107 down to 57-108 (dimension grows from A1:R107 to
# A1:R108), and populates the newly opened row 56 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 56..107 down by one to make room for the new record.
$ws.Rows.Item(56).Insert()

# Populate the new row 56 with the new observation.
$ws.Range("A56").Value = 2
$ws.Range("B56").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C56").Value = "Coquimbo"
$ws.Range("D56").Value = 45175
$ws.Range("E56").Value = 4
$ws.Range("F56").Value = 100112022
$ws.Range("G56").Value = "Arveja Verde"
$ws.Range("H56").Value = "Perfection"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 400
$ws.Range("K56").Value = 21000
$ws.Range("L56").Value = 23000
$ws.Range("M56").Value = 22000
$ws.Range("N56").Value = "$/malla 25 kilos"
$ws.Range("O56").Value = "Provincia de Limarí"
$ws.Range("P56").Value = 880
$ws.Range("Q56").Value = 25
$ws.Range("R56").Value = "Hortaliza"
